# New scenarios, data-output integration
# Adds a "BestChoice?" column (K) to the Connors test map metadata sheet,
# flagging every existing scenario row with option-3 as the best choice,
# now that data output is handled via a plugin call.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in K1
$ws.Range("K1").Value = "BestChoice?"

# Every data row (2..79) gets a BestChoice value of 3
for ($r = 2; $r -le 79; $r++) {
    $ws.Cells.Item($r, 11).Value = 3
}

# Resize the new column to fit its contents, like the other bestFit columns
$ws.Columns.Item(11).AutoFit() | Out-Null

# Move the view / active cell the way the author left it
$ws.Range("J3").Select() | Out-Null

Write-Output "BestChoice column added"
